# Update row 8 (ano/ano_obj = 2025) metrics on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 911
$ws.Range("D8").Value = 151
$ws.Range("E8").Value = 760
$ws.Range("F8").Value = 6.193601312551271
$ws.Range("G8").Value = 83.42480790340285
$ws.Range("H8").Value = 16.57519209659715
